$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix sort order of rows 6:7 (2004 Melody of Minstrel / 2003 Chains of Fortune were out of order) ---
$ws.Range("A6:F7").Sort($ws.Range("A6:A7"))

# --- Add new rows 31:41 (additional Blade of Arcana releases) ---
$ws.Range("A31").Value = 2006
$ws.Range("B31").Value = 'ブレイド・オブ・アルカナ The 3rd Edition'
$ws.Range("D31").Value = 'Enterbrain'
$ws.Range("E31").Value = '3rd_edition.jpg'
$ws.Range("C31").Value = 'Blade of Arcana: The 3rd Edition'
$ws.Range("A34").Value = 2009
$ws.Range("D34").Value = 'Game Field'
$ws.Range("E34").Value = 'gate_of_the_brave.jpg'
$ws.Range("B34").Value = 'ゲート・オブ・ザ・ブレイヴ'
$ws.Range("C34").Value = 'Gate of the Brave'
$ws.Range("A32").Value = 2008
$ws.Range("B32").Value = 'ランド・オブ・ザ・ギルティ'
$ws.Range("C32").Value = 'Land of the Guilty'
$ws.Range("D32").Value = 'Game Field'
$ws.Range("E32").Value = 'land_of_the_guilty_3rd.jpg'
$ws.Range("A33").Value = 2008
$ws.Range("D33").Value = 'Game Field'
$ws.Range("E33").Value = 'sun_of_darkness.jpg'
$ws.Range("B33").Value = 'サン・オブ・ダークネス'
$ws.Range("C33").Value = 'Sun of Darkness'
$ws.Range("A36").Value = 2010
$ws.Range("B36").Value = 'スレイヤーズ・オブ・レッドドラゴン'
$ws.Range("C36").Value = 'Slayers of Red Dragon'
$ws.Range("D36").Value = 'Enterbrain'
$ws.Range("E36").Value = 'slayers_of_red_dragon.jpg'
$ws.Range("A37").Value = 2012
$ws.Range("B37").Value = 'キング・オブ・ザ・ランド'
$ws.Range("C37").Value = 'King of the Land'
$ws.Range("D37").Value = 'Enterbrain'
$ws.Range("E37").Value = 'king_of_the_land.jpg'
$ws.Range("A35").Value = 2009
$ws.Range("B35").Value = ' 剣十字の騎士'
$ws.Range("C35").Value = 'Knight of the Sword Cross'
$ws.Range("D35").Value = 'Enterbrain'
$ws.Range("E35").Value = 'knight_of_the_sword_cross.jpg'
$ws.Range("F35").Value = 'replay'
$ws.Range("A40").Value = 2015
$ws.Range("F40").Value = 'replay'
$ws.Range("A38").Value = 2015
$ws.Range("E38").Value = 'blade_of_arcana_reincarnation.jpg'
$ws.Range("B38").Value = 'ブレイド・オブ・アルカナ リインカーネイション'
$ws.Range("C38").Value = 'Blade of Arcana Reincarnation'
$ws.Range("D38").Value = 'Kadokawa'
$ws.Range("D40").Value = 'Kadokawa'
$ws.Range("A41").Value = 2016
$ws.Range("D41").Value = 'Kadokawa'
$ws.Range("A39").Value = 2015
$ws.Range("E39").Value = 'ground_of_valor.jpg'
$ws.Range("B39").Value = 'グラウンド・オブ・ヴァラー'
$ws.Range("C39").Value = 'Ground of Valor'
$ws.Range("D39").Value = 'F.E.A.R.'
$ws.Range("E41").Value = 'crown_of_evil.jpg'
$ws.Range("B41").Value = 'クラウン・オブ・イビル'
$ws.Range("C41").Value = 'Crown of Evil'
$ws.Range("B40").Value = '刻まれし者の詩'
$ws.Range("C40").Value = 'Engraved Poetry'
$ws.Range("E40").Value = 'engraved_poetry.jpg'
$ws.Range("F31").Value = 'rulebook'
$ws.Range("F32").Value = 'supplement'
$ws.Range("F33").Value = 'supplement'
$ws.Range("F34").Value = 'supplement'
$ws.Range("F36").Value = 'supplement'
$ws.Range("F37").Value = 'supplement'
$ws.Range("F38").Value = 'rulebook'
$ws.Range("F39").Value = 'supplement'
$ws.Range("F41").Value = 'supplement'

# --- Sort the newly added rows (A31:F41) by year, ascending ---
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("A31:A41"))
$sortObj.SetRange($ws.Range("A31:F41"))
$sortObj.Header = -4142
$sortObj.Apply()

# --- Select the newly added rows, matching the final UI selection state ---
$ws.Range("A31:A41").EntireRow.Select()
